$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 13890768
$ws.Range("I33").Value = 22727758
$ws.Range("K33").Value = 22727758
$ws.Range("M33").Value = -22727529
$ws.Range("H46").Value = 1198.3334
$ws.Range("J46").Value = 1108
$ws.Range("L46").Value = 3324
$ws.Range("N46").Value = -3562
$ws.Range("H60").Value = 1198.3334
$ws.Range("J60").Value = 1108
$ws.Range("L60").Value = 3324
$ws.Range("N60").Value = -4292
$ws.Range("H62").Value = 3191.4167
$ws.Range("I62").Value = 2477.4443
$ws.Range("K62").Value = 2477.4443
$ws.Range("M62").Value = -1853.4443
$ws.Range("H64").Value = 5207.1665
$ws.Range("I64").Value = 2472
$ws.Range("K64").Value = 2472
$ws.Range("M64").Value = -2224
$ws.Range("H65").Value = 3191.4167
$ws.Range("I65").Value = 2477.4443
$ws.Range("K65").Value = 12387.2215
$ws.Range("M65").Value = -9267.2215
$ws.Range("H67").Value = 5207.1665
$ws.Range("I67").Value = 2472
$ws.Range("K67").Value = 2472
$ws.Range("M67").Value = -1614
$ws.Range("H86").Value = 7936.579
$ws.Range("J86").Value = 7380.3
$ws.Range("L86").Value = 7380.3
$ws.Range("N86").Value = -9626.299999999999
$ws.Range("H89").Value = 7936.579
$ws.Range("J89").Value = 7380.3
$ws.Range("L89").Value = 36901.5
$ws.Range("N89").Value = -48133.5
$ws.Range("H116").Value = 3882.9473
$ws.Range("I116").Value = 3495.8667
$ws.Range("K116").Value = 3495.8667
$ws.Range("M116").Value = -53.86670000000004
$ws.Range("H138").Value = 4971.32
$ws.Range("J138").Value = 4766.653
$ws.Range("L138").Value = 14299.959
$ws.Range("N138").Value = -24579.959

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4568716.5
$ws.Range("I132").Value = 1299.5
$ws.Range("K132").Value = 3898.5
$ws.Range("M132").Value = -1368.5

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 3748.5
$ws.Range("I3").Value = 3748.5
$ws.Range("K3").Value = 3748.5
$ws.Range("M3").Value = -3635.5
$ws.Range("H31").Value = 34576.5
$ws.Range("I31").Value = 1366
$ws.Range("J31").Value = 67787
$ws.Range("K31").Value = 1366
$ws.Range("L31").Value = 67787
$ws.Range("M31").Value = -1071
$ws.Range("N31").Value = -68377
$ws.Range("H34").Value = 34576.5
$ws.Range("I34").Value = 1366
$ws.Range("J34").Value = 67787
$ws.Range("K34").Value = 1366
$ws.Range("L34").Value = 67787
$ws.Range("M34").Value = -1164
$ws.Range("N34").Value = -68191
$ws.Range("H86").Value = 15513.286
$ws.Range("I86").Value = 16432.334
$ws.Range("J86").Value = 9999
$ws.Range("K86").Value = 16432.334
$ws.Range("L86").Value = 9999
$ws.Range("M86").Value = -15309.334
$ws.Range("N86").Value = -12245
$ws.Range("H89").Value = 15513.286
$ws.Range("I89").Value = 16432.334
$ws.Range("J89").Value = 9999
$ws.Range("K89").Value = 82161.67
$ws.Range("L89").Value = 49995
$ws.Range("M89").Value = -76545.67
$ws.Range("N89").Value = -61227
$ws.Range("H105").Value = 7844.875
$ws.Range("I105").Value = 13527.25
$ws.Range("J105").Value = 2162.5
$ws.Range("K105").Value = 13527.25
$ws.Range("L105").Value = 2162.5
$ws.Range("M105").Value = -11780.25
$ws.Range("N105").Value = -5656.5
$ws.Range("H107").Value = 1808.8148
$ws.Range("I107").Value = 1664.6471
$ws.Range("J107").Value = 2053.9
$ws.Range("K107").Value = 1664.6471
$ws.Range("L107").Value = 2053.9
$ws.Range("M107").Value = 255.3529000000001
$ws.Range("N107").Value = -5893.9

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1123.75
$ws.Range("J114").Value = 1123.75
$ws.Range("L114").Value = 3371.25
$ws.Range("N114").Value = -9879.25
$ws.Range("H129").Value = 10003000
$ws.Range("I129").Value = 16668334
$ws.Range("J129").Value = 4999.5
$ws.Range("K129").Value = 50005002
$ws.Range("L129").Value = 14998.5
$ws.Range("M129").Value = -50000002
$ws.Range("N129").Value = -24998.5
$ws.Range("H137").Value = 3288.1667
$ws.Range("J137").Value = 2893
$ws.Range("L137").Value = 8679
$ws.Range("N137").Value = -18879

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2428.8333
$ws.Range("I14").Value = 1432.125
$ws.Range("K14").Value = 1432.125
$ws.Range("M14").Value = -1264.125
$ws.Range("H20").Value = 739580.9
$ws.Range("I20").Value = 2506033.8
$ws.Range("J20").Value = 32999.7
$ws.Range("K20").Value = 2506033.8
$ws.Range("L20").Value = 32999.7
$ws.Range("M20").Value = -2505788.8
$ws.Range("N20").Value = -33489.7
$ws.Range("H23").Value = 5163.5
$ws.Range("J23").Value = 6268.5
$ws.Range("L23").Value = 6268.5
$ws.Range("N23").Value = -6714.5
$ws.Range("H48").Value = 35000
$ws.Range("J48").Value = 35000
$ws.Range("L48").Value = 35000
$ws.Range("N48").Value = -35970
$ws.Range("H49").Value = 21000
$ws.Range("J49").Value = 21000
$ws.Range("L49").Value = 21000
$ws.Range("N49").Value = -21368
$ws.Range("H70").Value = 8317.375
$ws.Range("I70").Value = 5206.4
$ws.Range("J70").Value = 13502.333
$ws.Range("K70").Value = 5206.4
$ws.Range("L70").Value = 13502.333
$ws.Range("M70").Value = -4936.4
$ws.Range("N70").Value = -14042.333
$ws.Range("H73").Value = 8317.375
$ws.Range("I73").Value = 5206.4
$ws.Range("J73").Value = 13502.333
$ws.Range("K73").Value = 5206.4
$ws.Range("L73").Value = 13502.333
$ws.Range("M73").Value = -4270.4
$ws.Range("N73").Value = -15374.333
$ws.Range("H126").Value = 7473.9546
$ws.Range("J126").Value = 6139.3335
$ws.Range("L126").Value = 18418.0005
$ws.Range("N126").Value = -23358.0005
$ws.Range("H132").Value = 642254.5600000001
$ws.Range("I132").Value = 6067.5
$ws.Range("K132").Value = 18202.5
$ws.Range("M132").Value = -15672.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 4303.8
$ws.Range("J21").Value = 3007
$ws.Range("L21").Value = 3007
$ws.Range("N21").Value = -3355
$ws.Range("H41").Value = 15000
$ws.Range("J41").Value = 15000
$ws.Range("L41").Value = 15000
$ws.Range("N41").Value = -15876
$ws.Range("H136").Value = 1148610.6
$ws.Range("I136").Value = 18249.385
$ws.Range("K136").Value = 54748.155
$ws.Range("M136").Value = -52198.155

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5020600
$ws.Range("I3").Value = 5020600
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5020600
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -5020486
$ws.Range("N3").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H132").Value = 579781.25
$ws.Range("I132").Value = 2748.6875
$ws.Range("K132").Value = 8246.0625
$ws.Range("M132").Value = -5716.0625
